# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Refresh COVID case counters for the affected countries
# - Argentina/Colombia and Bolivia/Chequia swap rank position (rows 9/10 and 36/37)
# - San Martin (Parte Francesa) moves above Comoras / Islas Feroe (rows 180-182)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 15 de Octubre de 2020 a las 01:19"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 8146377
$ws.Cells.Item(4, 3).Value = 56027
$ws.Cells.Item(4, 4).Value = 5269032
$ws.Cells.Item(4, 5).Value = 2655544
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 928
$ws.Cells.Item(4, 8).Value = 221801

# Row 6: Brasil
$ws.Cells.Item(6, 1).Value = "Brasil"
$ws.Cells.Item(6, 2).Value = 5141498
$ws.Cells.Item(6, 3).Value = 26675
$ws.Cells.Item(6, 4).Value = 4568813
$ws.Cells.Item(6, 5).Value = 420906
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 716
$ws.Cells.Item(6, 8).Value = 151779

# Row 9: Argentina
$ws.Cells.Item(9, 1).Value = "Argentina"
$ws.Cells.Item(9, 2).Value = 931967
$ws.Cells.Item(9, 3).Value = 14932
$ws.Cells.Item(9, 4).Value = 751146
$ws.Cells.Item(9, 5).Value = 155900
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 349
$ws.Cells.Item(9, 8).Value = 24921

# Row 10: Colombia
$ws.Cells.Item(10, 1).Value = "Colombia"
$ws.Cells.Item(10, 2).Value = 930159
$ws.Cells.Item(10, 3).Value = 6061
$ws.Cells.Item(10, 4).Value = 816667
$ws.Cells.Item(10, 5).Value = 85186
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 165
$ws.Cells.Item(10, 8).Value = 28306

# Row 36: Chequia
$ws.Cells.Item(36, 1).Value = "Chequia"
$ws.Cells.Item(36, 2).Value = 139290
$ws.Cells.Item(36, 3).Value = 9543
$ws.Cells.Item(36, 4).Value = 60901
$ws.Cells.Item(36, 5).Value = 77217
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 66
$ws.Cells.Item(36, 8).Value = 1172

# Row 37: Bolivia
$ws.Cells.Item(37, 1).Value = "Bolivia"
$ws.Cells.Item(37, 2).Value = 138922
$ws.Cells.Item(37, 3).Value = 227
$ws.Cells.Item(37, 4).Value = 102083
$ws.Cells.Item(37, 5).Value = 28488
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 25
$ws.Cells.Item(37, 8).Value = 8351

# Row 39: Panama
$ws.Cells.Item(39, 1).Value = "Panama"
$ws.Cells.Item(39, 2).Value = 122128
$ws.Cells.Item(39, 3).Value = 832
$ws.Cells.Item(39, 4).Value = 97919
$ws.Cells.Item(39, 5).Value = 21690
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 8
$ws.Cells.Item(39, 8).Value = 2519

# Row 46: Egipto
$ws.Cells.Item(46, 1).Value = "Egipto"
$ws.Cells.Item(46, 2).Value = 104915
$ws.Cells.Item(46, 3).Value = 128
$ws.Cells.Item(46, 4).Value = 97920
$ws.Cells.Item(46, 5).Value = 918
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = 6
$ws.Cells.Item(46, 8).Value = 6077

# Row 61: Nigeria
$ws.Cells.Item(61, 1).Value = "Nigeria"
$ws.Cells.Item(61, 2).Value = 60834
$ws.Cells.Item(61, 3).Value = 179
$ws.Cells.Item(61, 4).Value = 52143
$ws.Cells.Item(61, 5).Value = 7575
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 1116

# Row 91: Camerun
$ws.Cells.Item(91, 1).Value = "Camerun"
$ws.Cells.Item(91, 2).Value = 21441
$ws.Cells.Item(91, 3).Value = 238
$ws.Cells.Item(91, 4).Value = 20117
$ws.Cells.Item(91, 5).Value = 901
$ws.Cells.Item(91, 6).Value = 0
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = 423

# Row 96: Noruega
$ws.Cells.Item(96, 1).Value = "Noruega"
$ws.Cells.Item(96, 2).Value = 15953
$ws.Cells.Item(96, 3).Value = 162
$ws.Cells.Item(96, 4).Value = 11863
$ws.Cells.Item(96, 5).Value = 3813
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 277

# Row 106: Consejo Danes para los Refugiados
$ws.Cells.Item(106, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(106, 2).Value = 10935
$ws.Cells.Item(106, 3).Value = 63
$ws.Cells.Item(106, 4).Value = 10306
$ws.Cells.Item(106, 5).Value = 348
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 5
$ws.Cells.Item(106, 8).Value = 281

# Row 115: Zimbabue
$ws.Cells.Item(115, 1).Value = "Zimbabue"
$ws.Cells.Item(115, 2).Value = 8055
$ws.Cells.Item(115, 3).Value = 19
$ws.Cells.Item(115, 4).Value = 7640
$ws.Cells.Item(115, 5).Value = 184
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 1
$ws.Cells.Item(115, 8).Value = 231

# Row 136: Republica de Africa Central
$ws.Cells.Item(136, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(136, 2).Value = 4855
$ws.Cells.Item(136, 3).Value = 1
$ws.Cells.Item(136, 4).Value = 1924
$ws.Cells.Item(136, 5).Value = 2869
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 62

# Row 146: Guyana
$ws.Cells.Item(146, 1).Value = "Guyana"
$ws.Cells.Item(146, 2).Value = 3589
$ws.Cells.Item(146, 3).Value = 24
$ws.Cells.Item(146, 4).Value = 2487
$ws.Cells.Item(146, 5).Value = 996
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 106

# Row 156: Uruguay
$ws.Cells.Item(156, 1).Value = "Uruguay"
$ws.Cells.Item(156, 2).Value = 2388
$ws.Cells.Item(156, 3).Value = 51
$ws.Cells.Item(156, 4).Value = 2007
$ws.Cells.Item(156, 5).Value = 330
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 7).Value = 0
$ws.Cells.Item(156, 8).Value = 51

# Row 180: San Martin (Parte Francesa)
$ws.Cells.Item(180, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(180, 2).Value = 501
$ws.Cells.Item(180, 3).Value = 35
$ws.Cells.Item(180, 4).Value = 380
$ws.Cells.Item(180, 5).Value = 113
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 8

# Row 181: Comoras
$ws.Cells.Item(181, 1).Value = "Comoras"
$ws.Cells.Item(181, 2).Value = 496
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 484
$ws.Cells.Item(181, 5).Value = 5
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 7

# Row 182: Islas Feroe
$ws.Cells.Item(182, 1).Value = "Islas Feroe"
$ws.Cells.Item(182, 2).Value = 478
$ws.Cells.Item(182, 3).Value = 1
$ws.Cells.Item(182, 4).Value = 467
$ws.Cells.Item(182, 5).Value = 11
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 0

# Row 199: San Bartolome
$ws.Cells.Item(199, 1).Value = "San Bartolome"
$ws.Cells.Item(199, 2).Value = 67
$ws.Cells.Item(199, 3).Value = 2
$ws.Cells.Item(199, 4).Value = 55
$ws.Cells.Item(199, 5).Value = 12
$ws.Cells.Item(199, 6).Value = 0
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 0
